$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New shared-string text blocks used below (kept as variables for clarity).
# Order of first use below matches the shared-string append order needed.
# ---------------------------------------------------------------------------
$graphsHeader  = 'Graphs'
$lessonsHeader = 'Lessons'

$lessonThu = 'don''t estimate the online trading (after 2:30pm)'

$lessonWed = 'Good to stick to the trade plan, but sell all UGAZ and bought back DGAZ too early. The trend is obvious (from 8am -930am) but a small dip lure you to act too early. Make 15% a day is good, but you could have made 24%. Need more patience but some being conservative will save your ass in the long run'

$noteWed = 'A great uptrend for NG before the Thursday''s report, though the consensus for the report is around -200 but I think this heavily expected report is already priced in during today''s trade. Techinically it could reach Monday''s gap low at 3.661. Weather is warmer than normal which is not supportive, this could be a bull trap in short term. Hold until 1030 to buy DGAZ when it jump high. One thing to to notice is that, expiration day is coming close, holidays could surpress consumption and today''s DGAZ trade volumn is huge, well, UGAZ''s trade volum is also huge, something big is going on, be very careful! '

# ---------------------------------------------------------------------------
# New "Graphs" / "Lessons" header columns (G, F).
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = $graphsHeader
$ws.Range("G1").HorizontalAlignment = -4131
$ws.Range("G1").WrapText = $false

$ws.Range("F1").Value = $lessonsHeader
$ws.Range("F1").HorizontalAlignment = -4131
$ws.Range("F1").WrapText = $false

# Lesson learned, attached to the existing Tuesday/Wednesday row (row 5).
$ws.Range("F5").Value = $lessonThu
$ws.Range("F5").HorizontalAlignment = -4131
$ws.Range("F5").WrapText = $false

# ---------------------------------------------------------------------------
# New trade-plan row (row 3): Wednesday 12/21 -> Thursday report trade plan.
# ---------------------------------------------------------------------------
$ws.Range("F3").Value = $lessonWed
$ws.Range("F3").HorizontalAlignment = -4131
$ws.Range("F3").WrapText = $true

$ws.Range("E3").Value = $noteWed

$ws.Range("A3").Value = 20161221
$ws.Range("B3").Value = 20161221
$ws.Range("C3").Value = "Wednesday"
$ws.Range("D3").Value = "Thursday"
$ws.Range("G3").Value = 20161221

$ws.Rows.Item(3).RowHeight = 100.8

# ---------------------------------------------------------------------------
# Column sizing: E shrinks to make room for the new F (Lessons) column.
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 83.33
$ws.Columns.Item(6).ColumnWidth = 63.67

# ---------------------------------------------------------------------------
# View: drop the old top-left/active-cell pin, select E4 instead.
# ---------------------------------------------------------------------------
$ws.Range("E4").Select()

# ---------------------------------------------------------------------------
# Page setup: printable in portrait orientation.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
